# Apply updated cryptocurrency price/volume data (GitHub Actions scrape refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).NumberFormat = "@"
$ws.Cells.Item(2, 4).Value = '43.222.88'
$ws.Cells.Item(2, 5).Value = '  -1.19%  '
$ws.Cells.Item(3, 4).NumberFormat = "@"
$ws.Cells.Item(3, 4).Value = '2.276.15'
$ws.Cells.Item(3, 5).Value = '  -0.61%  '
$ws.Cells.Item(4, 5).Value = '  -0.40%  '
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = '111.65'
$ws.Cells.Item(5, 5).Value = '  -1.67%  '
$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = '263.73'
$ws.Cells.Item(6, 5).Value = '  -1.65%  '
$ws.Cells.Item(7, 4).NumberFormat = "@"
$ws.Cells.Item(7, 4).Value = '0.634'
$ws.Cells.Item(7, 5).Value = '  +1.80%  '
$ws.Cells.Item(8, 5).Value = '  +0.14%  '
$ws.Cells.Item(9, 4).NumberFormat = "@"
$ws.Cells.Item(9, 4).Value = '0.607'
$ws.Cells.Item(9, 5).Value = '  -2.21%  '
$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = '46.72'
$ws.Cells.Item(10, 5).Value = '  -2.50%  '
$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = '0.0935'
$ws.Cells.Item(11, 5).Value = '  -0.32%  '
$ws.Cells.Item(12, 5).Value = '  +4.58%  '
$ws.Cells.Item(13, 5).Value = '  +1.62%  '
$ws.Cells.Item(14, 4).NumberFormat = "@"
$ws.Cells.Item(14, 4).Value = '15.39'
$ws.Cells.Item(14, 5).Value = '  -1.41%  '
$ws.Cells.Item(15, 4).NumberFormat = "@"
$ws.Cells.Item(15, 4).Value = '2.619.17'
$ws.Cells.Item(15, 5).Value = '  -0.56%  '
$ws.Cells.Item(16, 4).NumberFormat = "@"
$ws.Cells.Item(16, 4).Value = '0.862'
$ws.Cells.Item(16, 5).Value = '  +2.00%  '
$ws.Cells.Item(17, 4).NumberFormat = "@"
$ws.Cells.Item(17, 4).Value = '2.269.16'
$ws.Cells.Item(17, 5).Value = '  -0.97%  '
$ws.Cells.Item(18, 4).NumberFormat = "@"
$ws.Cells.Item(18, 4).Value = '43.206.19'
$ws.Cells.Item(18, 5).Value = '  -0.91%  '
$ws.Cells.Item(19, 5).Value = '  -1.25%  '
$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = '6.73'
$ws.Cells.Item(20, 5).Value = '  +3.38%  '
$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = '71.81'
$ws.Cells.Item(21, 5).Value = '  -1.04%  '
$ws.Cells.Item(22, 5).Value = '  -0.72%  '
$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = '234.11'
$ws.Cells.Item(23, 5).Value = '  +0.61%  '
$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = '9.42'
$ws.Cells.Item(24, 5).Value = '  -3.23%  '
$ws.Cells.Item(25, 5).Value = '  +1.30%  '
$ws.Cells.Item(26, 5).Value = '  +1.89%  '
$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = '11.32'
$ws.Cells.Item(27, 5).Value = '  -2.70%  '
$ws.Cells.Item(28, 5).Value = '  -0.15%  '
$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = '40.46'
$ws.Cells.Item(29, 5).Value = '  -3.66%  '
$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).Value = '3.35'
$ws.Cells.Item(30, 5).Value = '  -1.56%  '
$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).Value = '2.24'
$ws.Cells.Item(31, 5).Value = '  -0.60%  '
$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = '172.91'
$ws.Cells.Item(32, 5).Value = '  -2.08%  '
$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).Value = '21.47'
$ws.Cells.Item(33, 5).Value = '  -0.39%  '
$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).Value = '0.0899'
$ws.Cells.Item(34, 5).Value = '  -3.11%  '
$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = '5.64'
$ws.Cells.Item(35, 5).Value = '  +1.48%  '
$ws.Cells.Item(36, 5).Value = '  +1.43%  '
$ws.Cells.Item(37, 2).Value = 'VeChain'
$ws.Cells.Item(37, 3).Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).Value = '0.0370'
$ws.Cells.Item(37, 5).Value = '  +3.52%  '
$ws.Cells.Item(38, 2).Value = 'RenderToken'
$ws.Cells.Item(38, 3).Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Cells.Item(38, 4).NumberFormat = "@"
$ws.Cells.Item(38, 4).Value = '4.62'
$ws.Cells.Item(38, 5).Value = '  -2.43%  '
$ws.Cells.Item(39, 5).Value = '  +4.88%  '
$ws.Cells.Item(40, 5).Value = '  -4.73%  '
$ws.Cells.Item(41, 5).Value = '  +6.92%  '
$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = '76.18'
$ws.Cells.Item(42, 5).Value = '  +5.17%  '
$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = '14.07'
$ws.Cells.Item(43, 5).Value = '  +1.32%  '
$ws.Cells.Item(44, 5).Value = '  -2.64%  '
$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = '6.07'
$ws.Cells.Item(45, 5).Value = '  -0.64%  '
$ws.Cells.Item(46, 5).Value = '  -0.12%  '
$ws.Cells.Item(47, 5).Value = '  -3.91%  '
$ws.Cells.Item(48, 2).Value = 'FraxShare'
$ws.Cells.Item(48, 3).Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = '8.54'
$ws.Cells.Item(48, 5).Value = '  -2.10%  '
$ws.Cells.Item(49, 2).Value = 'Aave'
$ws.Cells.Item(49, 3).Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = '102.26'
$ws.Cells.Item(49, 5).Value = '  -1.07%  '
$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = '1.26'
$ws.Cells.Item(50, 5).Value = '  +2.26%  '
$ws.Cells.Item(51, 5).Value = '  -0.74%  '
